# Banco_Dados_Qualidade_02-26.xlsx edit
# - Adds 4 new tracking columns (BD:BG) used to feed the new "Análise Refugo" page
# - Fixes OP column (H) on row 3 to be stored as a number instead of text
# - Appends 3 new inspection records (rows 4-6)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New header columns
# ---------------------------------------------------------------------------
$ws.Range("BD1").Value = "Código Castanha"
$ws.Range("BE1").Value = "Motivo_Usinagem"
$ws.Range("BF1").Value = "Motivo_Medida"
$ws.Range("BG1").Value = "Motivo_Outros"

# The header row (row 1) used to carry a bold / bordered / centered style.
# The refreshed sheet drops that formatting, so bring the header row back to
# the workbook's default (unstyled) look.
$ws.Range("A1:BG1").ClearFormats()

# ---------------------------------------------------------------------------
# Row 3 fix: OP (H) was stored as text "04711301003" - normalize to a number
# ---------------------------------------------------------------------------
$ws.Range("H3").Value = 4711301003

# ---------------------------------------------------------------------------
# New row 4
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "17/02/2026"
$ws.Range("C4").Value = "Alexandre"
$ws.Range("H4").Value = 4753001001
$ws.Range("I4").Value = "MECÂNICA INDUSTRIAL M.N. LTDA"
$ws.Range("N4").Value = 1
$ws.Range("Q4").Value = 0
$ws.Range("AB4").Value = " | "
$ws.Range("BD4").Value = 2025

# ---------------------------------------------------------------------------
# New row 5
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "17/02/2026"
$ws.Range("C5").Value = "Alexandre"
$ws.Range("H5").Value = 4758601001
$ws.Range("I5").Value = "VCI BRASIL INDUSTRIA E COMERCIO DE EMBALAGENS LTDA"
$ws.Range("N5").Value = 1
$ws.Range("Q5").Value = 0
$ws.Range("AB5").Value = " | "
$ws.Range("BD5").Value = 2025

# ---------------------------------------------------------------------------
# New row 6
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "18/02/2026"
$ws.Range("C6").Value = "Alexandre"
# OP keeps its leading zero -> force text storage (quote-prefix) instead of
# letting Excel coerce it to the number 4755001004.
$ws.Range("H6").Value = "'04755001004"
$ws.Range("I6").Value = "AUMAQRS LTDA"
$ws.Range("N6").Value = 1
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = "RETRABALHO"
$ws.Range("AB6").Value = "sddsdsd | dsdsdsd"
# Stored as text on this row (inconsistent data entry in the source sheet)
$ws.Range("BD6").Value = "'2025"
$ws.Range("BE6").Value = $false
$ws.Range("BF6").Value = $true
$ws.Range("BG6").Value = $false
